$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.478.81"
Set-TextValue $ws.Range("E2") "  -0.48%  "

Set-TextValue $ws.Range("D3") "1.827.22"
Set-TextValue $ws.Range("E3") "  -1.18%  "

Set-TextValue $ws.Range("E4") "  -0.56%  "

Set-TextValue $ws.Range("D5") "333.58"
Set-TextValue $ws.Range("E5") "  -0.07%  "

Set-TextValue $ws.Range("D6") "1.005"
Set-TextValue $ws.Range("E6") "  -0.51%  "

Set-TextValue $ws.Range("D7") "0.4581"
Set-TextValue $ws.Range("E7") "  -0.23%  "

Set-TextValue $ws.Range("D8") "0.3817"
Set-TextValue $ws.Range("E8") "  -1.94%  "

Set-TextValue $ws.Range("D9") "46.17"
Set-TextValue $ws.Range("E9") "  +0.37%  "

Set-TextValue $ws.Range("D10") "0.07816"
Set-TextValue $ws.Range("E10") "  -1.26%  "

Set-TextValue $ws.Range("E11") "  -4.24%  "

Set-TextValue $ws.Range("D12") "20.98"
Set-TextValue $ws.Range("E12") "  -2.75%  "

Set-TextValue $ws.Range("D13") "1.837.72"
Set-TextValue $ws.Range("E13") "  -0.30%  "

Set-TextValue $ws.Range("D14") "5.829"
Set-TextValue $ws.Range("E14") "  -2.05%  "

Set-TextValue $ws.Range("D15") "7.045"
Set-TextValue $ws.Range("E15") "  -1.91%  "

Set-TextValue $ws.Range("D16") "1.007"
Set-TextValue $ws.Range("E16") "  -0.58%  "

Set-TextValue $ws.Range("E17") "  +1.05%  "

Set-TextValue $ws.Range("D18") "0.06588"
Set-TextValue $ws.Range("E18") "  -1.90%  "

Set-TextValue $ws.Range("D19") "0.00001018"
Set-TextValue $ws.Range("E19") "  -1.74%  "

Set-TextValue $ws.Range("D20") "17.06"
Set-TextValue $ws.Range("E20") "  -0.52%  "

Set-TextValue $ws.Range("D21") "1.005"
Set-TextValue $ws.Range("E21") "  -0.70%  "

Set-TextValue $ws.Range("D22") "27.469.10"
Set-TextValue $ws.Range("E22") "  -0.47%  "

Set-TextValue $ws.Range("D23") "5.280"
Set-TextValue $ws.Range("E23") "  -2.44%  "

Set-TextValue $ws.Range("D24") "10.76"
Set-TextValue $ws.Range("E24") "  -1.27%  "

Set-TextValue $ws.Range("D25") "2.273"
Set-TextValue $ws.Range("E25") "  -1.33%  "

Set-TextValue $ws.Range("D26") "2.072.35"
Set-TextValue $ws.Range("E26") "  +0.36%  "

Set-TextValue $ws.Range("D27") "158.84"
Set-TextValue $ws.Range("E27") "  -0.34%  "

Set-TextValue $ws.Range("E28") "  -1.26%  "

Set-TextValue $ws.Range("D29") "2.029"
Set-TextValue $ws.Range("E29") "  -4.53%  "

Set-TextValue $ws.Range("D30") "5.269"
Set-TextValue $ws.Range("E30") "  -3.02%  "

Set-TextValue $ws.Range("D31") "117.63"
Set-TextValue $ws.Range("E31") "  -3.21%  "

Set-TextValue $ws.Range("D32") "0.09377"
Set-TextValue $ws.Range("E32") "  -0.14%  "

Set-TextValue $ws.Range("D33") "0.9256"
Set-TextValue $ws.Range("E33") "  -4.85%  "

Set-TextValue $ws.Range("D34") "3.578"
Set-TextValue $ws.Range("E34") "  -1.16%  "

Set-TextValue $ws.Range("D35") "5.187"
Set-TextValue $ws.Range("E35") "  -2.11%  "

Set-TextValue $ws.Range("D36") "1.310"
Set-TextValue $ws.Range("E36") "  -1.74%  "

Set-TextValue $ws.Range("D37") "0.05946"
Set-TextValue $ws.Range("E37") "  -0.83%  "

Set-TextValue $ws.Range("D38") "0.02173"
Set-TextValue $ws.Range("E38") "  -2.34%  "

Set-TextValue $ws.Range("D39") "8.102"
Set-TextValue $ws.Range("E39") "  -3.02%  "

Set-TextValue $ws.Range("D40") "1.005"
Set-TextValue $ws.Range("E40") "  -0.59%  "

Set-TextValue $ws.Range("D41") "1.141"
Set-TextValue $ws.Range("E41") "  -4.19%  "

Set-TextValue $ws.Range("D42") "0.5709"
Set-TextValue $ws.Range("E42") "  -3.29%  "

Set-TextValue $ws.Range("E43") "  -2.27%  "

Set-TextValue $ws.Range("D44") "9.874"
Set-TextValue $ws.Range("E44") "  -5.27%  "

Set-TextValue $ws.Range("E45") "  +2.06%  "

Set-TextValue $ws.Range("D46") "0.5370"
Set-TextValue $ws.Range("E46") "  -3.58%  "

Set-TextValue $ws.Range("D47") "11.77"
Set-TextValue $ws.Range("E47") "  -2.64%  "

Set-TextValue $ws.Range("E48") "  -0.75%  "

Set-TextValue $ws.Range("D49") "0.06846"
Set-TextValue $ws.Range("E49") "  +2.11%  "

Set-TextValue $ws.Range("D50") "109.88"
Set-TextValue $ws.Range("E50") "  -1.07%  "

Set-TextValue $ws.Range("D51") "1.005"
Set-TextValue $ws.Range("E51") "  -32.60%  "
